# Apply the edits described by the diff:
#  - Metadata sheet: update the Date value
#  - Concepts sheet: fix three Display text typos

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Date (Metadata!B8): 2022-12-09T21:34:05+00:00 -> 2023-01-23T16:13:19+00:00
$wsMetadata.Range("B8").Value = "2023-01-23T16:13:19+00:00"

# Concepts!C3 (EU): "European Caucasia" -> "European Caucasian"
$wsConcepts.Range("C3").Value = "European Caucasian"

# Concepts!C4 (AFR): "African or Carabean" -> "African or Caribbean"
$wsConcepts.Range("C4").Value = "African or Caribbean"

# Concepts!C5 (LAT-AM): "Hispanic and Latino Americans" -> "Hispanic and Latino American"
$wsConcepts.Range("C5").Value = "Hispanic and Latino American"
